$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6802170276641846
$ws.Range("B1").Value = 0.4580024778842926
$ws.Range("C1").Value = 0.3558070659637451
$ws.Range("D1").Value = 0.3455876111984253
$ws.Range("E1").Value = 0.3739714324474335
